# Update the cryptos price list (columns B-E, rows 2-51) with refreshed
# values from the GitHub Actions data pull.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'" + '58.681.18'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -3.70%  '
$ws.Range('D3').Value = "'" + '2.556.50'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -1.77%  '
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('D5').Value = "'" + '504.73'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -3.66%  '
$ws.Range('D6').Value = "'" + '142.43'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -8.20%  '
$ws.Range('D7').Value = "'" + '0.999'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.10%  '
$ws.Range('D8').Value = "'" + '0.552'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -6.20%  '
$ws.Range('D9').Value = "'" + '2.563.55'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -1.62%  '
$ws.Range('D10').Value = "'" + '6.20'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -7.23%  '
$ws.Range('D11').Value = "'" + '0.101'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -4.38%  '
$ws.Range('D12').Value = "'" + '0.329'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -5.32%  '
$ws.Range('E13').Value = '  -1.09%  '
$ws.Range('D14').Value = "'" + '3.010.63'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -1.54%  '
$ws.Range('D15').Value = "'" + '58.696.36'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -3.67%  '
$ws.Range('D16').Value = "'" + '20.50'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -5.61%  '
$ws.Range('D17').Value = "'" + '0.0000134'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -5.42%  '
$ws.Range('D18').Value = "'" + '2.561.32'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -1.65%  '
$ws.Range('D19').Value = "'" + '4.50'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -5.82%  '
$ws.Range('D20').Value = "'" + '331.22'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -6.80%  '
$ws.Range('D21').Value = "'" + '10.02'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -5.47%  '
$ws.Range('D22').Value = "'" + '0.999'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.05%  '
$ws.Range('D23').Value = "'" + '5.92'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -4.85%  '
$ws.Range('D24').Value = "'" + '59.31'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -2.86%  '
$ws.Range('D25').Value = "'" + '0.405'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -4.89%  '
$ws.Range('D26').Value = "'" + '1.00'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.14%  '
$ws.Range('E27').Value = '  -5.75%  '
$ws.Range('D28').Value = "'" + '0.0₃0774'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -8.77%  '
$ws.Range('D29').Value = "'" + '6.84'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -7.73%  '
$ws.Range('D30').Value = "'" + '1.00'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -0.03%  '
$ws.Range('B31').Value = 'Monero'
$ws.Range('C31').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D31').Value = "'" + '149.44'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.37%  '
$ws.Range('B32').Value = 'EthereumClassic'
$ws.Range('C32').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D32').Value = "'" + '18.53'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -4.68%  '
$ws.Range('D33').Value = "'" + '5.81'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -7.36%  '
$ws.Range('D34').Value = "'" + '1.53'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -4.37%  '
$ws.Range('D35').Value = "'" + '3.92'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -7.00%  '
$ws.Range('D36').Value = "'" + '0.881'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -3.64%  '
$ws.Range('D37').Value = "'" + '1.10'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -8.69%  '
$ws.Range('D38').Value = "'" + '35.78'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -1.95%  '
$ws.Range('D39').Value = "'" + '0.821'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -8.54%  '
$ws.Range('D40').Value = "'" + '285.14'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -2.16%  '
$ws.Range('D41').Value = "'" + '1.38'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -8.59%  '
$ws.Range('D42').Value = "'" + '3.48'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -8.76%  '
$ws.Range('B43').Value = 'FirstDigitalUSD'
$ws.Range('C43').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D43').Value = "'" + '1.00'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.21%  '
$ws.Range('B44').Value = 'Mantle'
$ws.Range('C44').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D44').Value = "'" + '0.607'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -2.47%  '
$ws.Range('D45').Value = "'" + '0.0977'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -3.64%  '
$ws.Range('D46').Value = "'" + '0.0529'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -5.60%  '
$ws.Range('B47').Value = 'WhiteBITCoin'
$ws.Range('C47').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D47').Value = "'" + '10.34'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.08%  '
$ws.Range('B48').Value = 'EnergySwap'
$ws.Range('C48').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D48').Value = "'" + '18.61'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -4.94%  '
$ws.Range('D49').Value = "'" + '0.0225'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -5.35%  '
$ws.Range('D50').Value = "'" + '4.51'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -8.54%  '
$ws.Range('D51').Value = "'" + '1.907.95'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -2.73%  '
